$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign the ink colour rows (3-6) so their Material/Inventoryitem values
# appear in the new order: Black, Cyan, Yellow, Magenta.
$ws.Range("D3").Value2 = "Black - UV - "
$ws.Range("H3").Value2 = "10001817 - 9409 MIXING BLACK UV - INK"

$ws.Range("D4").Value2 = "Cyan - UV - "
$ws.Range("H4").Value2 = "10001837 - 9443 PRO CYAN BW8 UV - INK"

$ws.Range("D5").Value2 = "Yellow - UV - "
$ws.Range("H5").Value2 = "10001305 - PROCESS YELLOW C UV"

$ws.Range("D6").Value2 = "Magenta - UV - "
$ws.Range("H6").Value2 = "10001836 - 9442 PRO MAGENTA BW5 UV - INK"

# Reassign the glue rows (7-8) so Adhesive now comes before Silicone.
$ws.Range("D7").Value2 = "Adhesive"
$ws.Range("H7").Value2 = "10001053 - RAVENWOOD LINERLESS-ADHESIVE - 7445HD"

$ws.Range("D8").Value2 = "Silicone"
$ws.Range("H8").Value2 = "10016451 - Evonik RW 10 Teco RC Silicone"
